$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
    $ws.Cells.Item($r, 10).Value = $r + 7
}

for ($c = 1; $c -le 10; $c++) {
    $ws.Cells.Item(19, $c).Value = 19 + ($c - 1)
}

$ws.Range("L16").Select()
